$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$data = @(
    @{ row=2; a="8669759"; m="english"; n="DOGS" },
    @{ row=3; a="8df405e"; m="japanese"; n="FISH" },
    @{ row=4; a="887529c"; m="english"; n="REPTILES" },
    @{ row=5; a="bad2608"; m="japanese"; n="FISH" },
    @{ row=6; a="897dada"; m="english"; n="FISH" },
    @{ row=7; a="3fa3da7"; m="english"; n="DOGS" },
    @{ row=8; a="d485f6c"; m="english"; n="BIRDS" },
    @{ row=9; a="23298ef"; m="japanese"; n="FISH" },
    @{ row=10; a="0087ae9"; m="japanese"; n="REPTILES" },
    @{ row=11; a="19c9f17"; m="english"; n="FISH" },
    @{ row=12; a="b4f83df"; m="japanese"; n="REPTILES" },
    @{ row=13; a="6d76ef1"; m="english"; n="BIRDS" },
    @{ row=14; a="86158a5"; m="english"; n="BIRDS" },
    @{ row=15; a="d73fb62"; m="english"; n="FISH" },
    @{ row=16; a="7e3f6df"; m="english"; n="REPTILES" },
    @{ row=17; a="0c59f17"; m="english"; n="DOGS" },
    @{ row=18; a="fdf24c6"; m="japanese"; n="DOGS" },
    @{ row=19; a="7aca76d"; m="japanese"; n="FISH" },
    @{ row=20; a="b838df9"; m="japanese"; n="BIRDS" },
    @{ row=21; a="1774f01"; m="japanese"; n="FISH" },
    @{ row=22; a="a7c8e7f"; m="japanese"; n="REPTILES" },
    @{ row=23; a="030b13d"; m="english"; n="DOGS" },
    @{ row=24; a="d24b22e"; m="english"; n="REPTILES" },
    @{ row=25; a="1d03e19"; m="japanese"; n="FISH" },
    @{ row=26; a="7a0c892"; m="japanese"; n="BIRDS" },
    @{ row=27; a="10ca2d1"; m="japanese"; n="DOGS" },
    @{ row=28; a="0e7c367"; m="japanese"; n="BIRDS" },
    @{ row=29; a="ce8b644"; m="english"; n="DOGS" },
    @{ row=30; a="35670ba"; m="japanese"; n="REPTILES" },
    @{ row=31; a="fcbc79c"; m="english"; n="DOGS" },
    @{ row=32; a="f759528"; m="japanese"; n="DOGS" },
    @{ row=33; a="a18e571"; m="japanese"; n="FISH" },
    @{ row=34; a="cf25f01"; m="english"; n="DOGS" },
    @{ row=35; a="b6cf6a3"; m="english"; n="CATS" },
    @{ row=36; a="875fb6d"; m="japanese"; n="CATS" },
    @{ row=37; a="e339cad"; m="english"; n="DOGS" },
    @{ row=38; a="6245ab9"; m="japanese"; n="DOGS" },
    @{ row=39; a="b211a4b"; m="english"; n="REPTILES" },
    @{ row=40; a="89dd132"; m="japanese"; n="FISH" },
    @{ row=41; a="802cbea"; m="japanese"; n="DOGS" },
    @{ row=42; a="6f45921"; m="japanese"; n="CATS" },
    @{ row=43; a="2e0b118"; m="japanese"; n="DOGS" },
    @{ row=44; a="5795e67"; m="japanese"; n="DOGS" },
    @{ row=45; a="4d61391"; m="english"; n="FISH" },
    @{ row=46; a="e039713"; m="japanese"; n="BIRDS" },
    @{ row=47; a="8b143f3"; m="english"; n="BIRDS" },
    @{ row=48; a="4beb29d"; m="english"; n="BIRDS" },
    @{ row=49; a="a5e45ce"; m="japanese"; n="REPTILES" },
    @{ row=50; a="6bb9d5b"; m="english"; n="DOGS" },
    @{ row=51; a="e0f5542"; m="english"; n="DOGS" }
)

# Some new column-A hash values (e.g. "8669759", "5795e67") look like plain numbers
# or scientific notation to Excel's input parser, which would silently convert them
# to numeric values and strip their text typing/formatting. Force those specific cells
# to stay text, then restore their original style by pasting formats from an
# untouched donor cell that already carries the same style.
$numericFix = @{
    2 = "8669759"
    44 = "5795e67"
}

foreach ($r in $numericFix.Keys) {
    $donor = $ws.Cells.Item(3, 1)
    $target = $ws.Cells.Item($r, 1)
    $target.NumberFormat = "@"
    $target.Value2 = $numericFix[$r]
    $donor.Copy()
    $target.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

foreach ($item in $data) {
    if ($numericFix.ContainsKey($item.row) -eq $false) {
        $ws.Cells.Item($item.row, 1).Value2 = $item.a
    }
    $ws.Cells.Item($item.row, 13).Value2 = $item.m
    $ws.Cells.Item($item.row, 14).Value2 = $item.n
}

Write-Output "done"